$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue $ws "D2" "23.907.65"
Set-TextValue $ws "E2" "  -2.03%  "
Set-TextValue $ws "D3" "1.651.35"
Set-TextValue $ws "E3" "  -0.89%  "
Set-TextValue $ws "D4" "1.002"
Set-TextValue $ws "E4" "  +0.05%  "
Set-TextValue $ws "D5" "310.12"
Set-TextValue $ws "E5" "  -0.97%  "
Set-TextValue $ws "D6" "1.002"
Set-TextValue $ws "E6" "  +0.06%  "
Set-TextValue $ws "D7" "0.3881"
Set-TextValue $ws "E7" "  -1.89%  "
Set-TextValue $ws "D8" "0.3815"
Set-TextValue $ws "E8" "  -2.34%  "
Set-TextValue $ws "D9" "51.95"
Set-TextValue $ws "E9" "  -0.77%  "
Set-TextValue $ws "D10" "1.350"
Set-TextValue $ws "E10" "  -3.50%  "
Set-TextValue $ws "D11" "1.002"
Set-TextValue $ws "E11" "  +0.08%  "
Set-TextValue $ws "D12" "0.08460"
Set-TextValue $ws "E12" "  -1.44%  "
Set-TextValue $ws "D13" "23.93"
Set-TextValue $ws "E13" "  -1.81%  "
Set-TextValue $ws "D14" "7.076"
Set-TextValue $ws "E14" "  -3.06%  "
Set-TextValue $ws "D15" "8.043"
Set-TextValue $ws "E15" "  +1.18%  "
Set-TextValue $ws "D16" "0.00001313"
Set-TextValue $ws "E16" "  -2.14%  "
Set-TextValue $ws "D17" "1.647.80"
Set-TextValue $ws "E17" "  -0.94%  "
Set-TextValue $ws "D18" "94.24"
Set-TextValue $ws "E18" "  -1.04%  "
Set-TextValue $ws "D19" "0.06996"
Set-TextValue $ws "E19" "  -0.14%  "
Set-TextValue $ws "D20" "19.65"
Set-TextValue $ws "E20" "  -4.36%  "
Set-TextValue $ws "D21" "6.952"
Set-TextValue $ws "E21" "  -0.58%  "
Set-TextValue $ws "D23" "13.76"
Set-TextValue $ws "E23" "  -0.02%  "
Set-TextValue $ws "D24" "23.898.82"
Set-TextValue $ws "E24" "  -2.08%  "
Set-TextValue $ws "D25" "2.447"
Set-TextValue $ws "E25" "  +1.20%  "
Set-TextValue $ws "D26" "2.968"
Set-TextValue $ws "E26" "  -2.79%  "
Set-TextValue $ws "D27" "22.09"
Set-TextValue $ws "E27" "  -2.08%  "
Set-TextValue $ws "D28" "153.51"
Set-TextValue $ws "E28" "  -2.29%  "
Set-TextValue $ws "D29" "5.415"
Set-TextValue $ws "E29" "  +0.00%  "
Set-TextValue $ws "D30" "138.01"
Set-TextValue $ws "E30" "  -3.23%  "
Set-TextValue $ws "D31" "7.858"
Set-TextValue $ws "E31" "  -2.35%  "
Set-TextValue $ws "D32" "2.506"
Set-TextValue $ws "E32" "  -1.56%  "
Set-TextValue $ws "D33" "1.836.68"
Set-TextValue $ws "E33" "  -0.48%  "
Set-TextValue $ws "D34" "1.018"
Set-TextValue $ws "E34" "  -3.49%  "
Set-TextValue $ws "D35" "0.08174"
Set-TextValue $ws "E35" "  -1.00%  "
$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws "D36" "6.726"
Set-TextValue $ws "E36" "  -2.88%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws "D37" "0.02927"
Set-TextValue $ws "E37" "  -3.19%  "
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws "D38" "10.84"
Set-TextValue $ws "E38" "  -3.00%  "
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws "D39" "0.2680"
Set-TextValue $ws "E39" "  -2.83%  "
Set-TextValue $ws "D40" "0.09122"
Set-TextValue $ws "E40" "  -1.10%  "
Set-TextValue $ws "D41" "0.7583"
Set-TextValue $ws "E41" "  -1.78%  "
Set-TextValue $ws "D42" "13.48"
Set-TextValue $ws "E42" "  -2.21%  "
Set-TextValue $ws "D43" "1.423"
Set-TextValue $ws "E43" "  -1.62%  "
Set-TextValue $ws "D44" "16.24"
Set-TextValue $ws "E44" "  -1.57%  "
Set-TextValue $ws "D45" "0.6943"
Set-TextValue $ws "E45" "  -2.46%  "
Set-TextValue $ws "D46" "2.463"
Set-TextValue $ws "E46" "  -3.03%  "
Set-TextValue $ws "D47" "4.093"
Set-TextValue $ws "E47" "  -0.98%  "
Set-TextValue $ws "D48" "1.001"
Set-TextValue $ws "E48" "  +0.03%  "
Set-TextValue $ws "D49" "0.08298"
Set-TextValue $ws "E49" "  -1.49%  "
Set-TextValue $ws "D50" "134.50"
Set-TextValue $ws "E50" "  -1.49%  "
Set-TextValue $ws "D51" "1.229"
Set-TextValue $ws "E51" "  -2.90%  "
